# Support for multiline headers
# Insert a new "Multi Line" worksheet between Sheet1 and Sheet2, with a
# two-row (multi-line) header, and make it the active sheet.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")
$newSheet = $wb.Worksheets.Add($null, $sheet1)
$newSheet.Name = "Multi Line"

# Header row, split across two rows so columns A/B read as multi-line
# headers ("First" + "Part" = "First Part", "Second" + "Part" = "Second
# Part") while column C's header is a single line ("Third Part").
$newSheet.Range("A1").Value = "First"
$newSheet.Range("A2").Value = "Part"
$newSheet.Range("B1").Value = "Second"
$newSheet.Range("B2").Value = "Part"
$newSheet.Range("C1").Value = "Third Part"

# Data rows
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "A"
$newSheet.Range("C3").Value = "X"

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "B"
$newSheet.Range("C4").Value = "Y"

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "C"
$newSheet.Range("C5").Value = "Z"

# Match the saved selection/active-cell on the new sheet.
[void]$newSheet.Range("C5").Select()
